$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.050.74"
$ws.Range("E2").Value = "  -0.89%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.195.24"
$ws.Range("E3").Value = "  -0.01%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.43"
$ws.Range("E5").Value = "  +1.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.87"
$ws.Range("E6").Value = "  +0.43%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.189.26"
$ws.Range("E8").Value = "  +0.21%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.551"
$ws.Range("E9").Value = "  -1.06%  "
$ws.Range("E10").Value = "  -0.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.63"
$ws.Range("E11").Value = "  -5.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.502"
$ws.Range("E12").Value = "  -3.33%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000269"
$ws.Range("E13").Value = "  -0.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.26"
$ws.Range("E14").Value = "  -2.90%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.724.15"
$ws.Range("E15").Value = "  +0.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.284.64"
$ws.Range("E16").Value = "  -0.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.33"
$ws.Range("E17").Value = "  -2.49%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.201.42"
$ws.Range("E18").Value = "  +0.06%  "
$ws.Range("E19").Value = "  +1.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "506.23"
$ws.Range("E20").Value = "  -2.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.24"
$ws.Range("E21").Value = "  -1.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.729"
$ws.Range("E22").Value = "  -1.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.97"
$ws.Range("E23").Value = "  -1.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.55"
$ws.Range("E24").Value = "  -3.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.89"
$ws.Range("E25").Value = "  -1.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.99"
$ws.Range("E27").Value = "  -0.67%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.99"
$ws.Range("E28").Value = "  -3.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.34"
$ws.Range("E29").Value = "  -0.90%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.126"
$ws.Range("E30").Value = "  +39.64%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.91"
$ws.Range("E31").Value = "  -0.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.95"
$ws.Range("E32").Value = "  -2.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "28.08"
$ws.Range("E33").Value = "  -1.13%  "
$ws.Range("E34").Value = "  +0.20%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.17"
$ws.Range("E35").Value = "  -5.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.46"
$ws.Range("E36").Value = "  -1.43%  "
$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "55.29"
$ws.Range("E37").Value = "  +0.49%  "
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "497.77"
$ws.Range("E38").Value = "  -3.84%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0764"
$ws.Range("E39").Value = "  +12.95%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0420"
$ws.Range("E40").Value = "  -1.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.130"
$ws.Range("E41").Value = "  +0.56%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.01"
$ws.Range("E42").Value = "  +2.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.70"
$ws.Range("E43").Value = "  -2.65%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.297"
$ws.Range("E44").Value = "  -2.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.909.93"
$ws.Range("E45").Value = "  -0.17%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.42"
$ws.Range("E46").Value = "  -1.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.08"
$ws.Range("E47").Value = "  -2.66%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.40"
$ws.Range("E48").Value = "  +1.70%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.116"
$ws.Range("E50").Value = "  -1.48%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "122.21"
$ws.Range("E51").Value = "  -0.35%  "
